# feat: add 2022-Q1 data
#
# - Old "总计" sheet (aggregate quarterly table) is repurposed into a new
#   "2022-Q1" sheet holding the per-fund holdings detail for 2022-Q1
#   (same layout as the other quarterly sheets, e.g. "2021-Q4").
# - A brand-new "总计" sheet is appended at the end, containing the same
#   aggregate table as before plus a new top row for "2022-Q1".

$wb = $excel.ActiveWorkbook

# Use the untouched "2021-Q4" sheet as a style template: its header row
# (B1:H1) and its "index" column (A2:A4) both carry the workbook's
# standard bold/centered/bordered style used throughout every sheet.
$template = $wb.Worksheets.Item("2021-Q4")

# The existing "总计" sheet will be fully rebuilt (it is recreated below),
# so just drop it now; its old content is re-entered verbatim further down.
$oldTotal = $wb.Worksheets.Item("总计")
[void]$oldTotal.Delete()

# ---------------------------------------------------------------------
# 1) New sheet "2022-Q1" — per-fund holdings detail, inserted right after
#    "2021-Q4" (i.e. at the end of the quarter sheets, before "总计").
# ---------------------------------------------------------------------
$lastQuarterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastQuarterSheet)
$q1.Name = "2022-Q1"

# Header row style (B1:H1) + index column style (A2:A6), copied from the
# template sheet so the new sheet matches the look of its siblings.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$template.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Helper-free: data rows (A = 0-based index, numeric; B..G are text-typed
# exactly like in the other quarter sheets; H is numeric).
$q1Rows = @(
    @{ A = 0; B = "002423"; C = "华宝兴业标普美国消费(QDII-LOF)美元"; D = "3.62"; E = "94.37"; F = "2.20"; G = "0.0796"; H = 9 },
    @{ A = 1; B = "162415"; C = "华宝标普美国消费(QDII-LOF)人民币A"; D = "3.62"; E = "94.37"; F = "2.20"; G = "0.0796"; H = 9 },
    @{ A = 2; B = "009975"; C = "华宝标普美国消费(QDII-LOF)人民币C"; D = "0.61"; E = "94.37"; F = "2.20"; G = "0.0134"; H = 9 },
    @{ A = 3; B = "519981"; C = "长信美国标准普尔100等权重指数增强(QDII)"; D = "0.47"; E = "84.16"; F = "0.88"; G = "0.0041"; H = 5 },
    @{ A = 4; B = "011706"; C = "长信美国标准普尔100等权重指数增强(QDII) - 美元"; D = "0.47"; E = "84.16"; F = "0.88"; G = "0.0041"; H = 5 }
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row.A

    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row.B
    $q1.Cells.Item($r, 2).Style = "Normal"

    $q1.Cells.Item($r, 3).Value = $row.C

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row.D
    $q1.Cells.Item($r, 4).Style = "Normal"

    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row.E
    $q1.Cells.Item($r, 5).Style = "Normal"

    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row.F
    $q1.Cells.Item($r, 6).Style = "Normal"

    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row.G
    $q1.Cells.Item($r, 7).Style = "Normal"

    $q1.Cells.Item($r, 8).Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) New sheet "总计" — aggregate table, appended at the very end.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$total.Name = "总计"

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$template.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @{ A = 0; B = "2022-Q1"; C = 5; D = 0.18 },
    @{ A = 1; B = "2021-Q4"; C = 3; D = 0.19 },
    @{ A = 2; B = "2021-Q3"; C = 3; D = 0.17 },
    @{ A = 3; B = "2021-Q2"; C = 3; D = 0.14 },
    @{ A = 4; B = "2021-Q1"; C = 3; D = 0.15 },
    @{ A = 5; B = "2020-Q4"; C = 4; D = 0.13 }
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row.A
    $total.Cells.Item($r, 2).Value = $row.B
    $total.Cells.Item($r, 3).Value = $row.C
    $total.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}
